# Spilamberto.xlsx - "aggiornamento fino a 20/09/2021"
# Append rows 375-385 (new daily data points) to the end of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data: row, date-serial (col A), nuovi pos. (col B), somma mobile 7gg. (col C),
# somma mobile 7gg. per 100mila abitanti (col D)
$data = @(
    @(375, 44449, 3, 18, 141.3427561837456),
    @(376, 44450, 0, 15, 117.7856301531213),
    @(377, 44451, 0, 14, 109.9332548095799),
    @(378, 44452, 1, 9, 70.67137809187278),
    @(379, 44453, 3, 7, 54.96662740478995),
    @(380, 44454, 0, 7, 54.96662740478995),
    @(381, 44455, 2, 9, 70.67137809187278),
    @(382, 44456, 6, 12, 94.22850412249706),
    @(383, 44457, 5, 17, 133.4903808402042),
    @(384, 44458, 0, 17, 133.4903808402042),
    @(385, 44459, 2, 18, 141.3427561837456)
)

$firstRow = $data[0][0]
$lastRow = $data[$data.Count - 1][0]

# Carry the existing date-column formatting (style "s=2": bordered, centered,
# bold, custom date/time number format) down onto the new A-column cells by
# copying the format of the last existing data row before writing values.
$ws.Range("A374").Copy($ws.Range(("A{0}:A{1}" -f $firstRow, $lastRow)))

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}
